# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Rebuild the worker/period data block (rows 16-33) on Hoja1: Luz Marina Torres
# Angarita's rows now come first (periods 2105..2010, descending), followed by
# Guillermo Jesus Rivero Gamarra's rows (periods 2106..2010, descending).  Also
# corrects the "Valor Mora" amounts for Luz Marina's 2105/2106 periods, which
# had been swapped (2105 was 35112 / 2106 was 24578; now 2105 is 24578 and
# 2106 is 35112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, DocId (col C), Name (col D), Period (col E), Valor Mora (col F), Salario Basico (col G)
$data = @(
    @(16, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2105", 24578, 877803),
    @(17, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2106", 35112, 877803),
    @(18, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2104", 35112, 877803),
    @(19, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2103", 35112, 877803),
    @(20, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2102", 35112, 877803),
    @(21, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2101", 35112, 877803),
    @(22, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2012", 35112, 877803),
    @(23, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2011", 35112, 877803),
    @(24, "45691962",   "LUZ MARINA TORRES ANGARITA",     "2010", 35112, 877803),
    @(25, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2106", 28000, 1000000),
    @(26, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2105", 40000, 1000000),
    @(27, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2104", 40000, 1000000),
    @(28, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2103", 40000, 1000000),
    @(29, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2102", 40000, 1000000),
    @(30, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2101", 40000, 1000000),
    @(31, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2012", 40000, 1000000),
    @(32, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2011", 40000, 1000000),
    @(33, "1047470977", "GUILLERMO JESUS RIVERO GAMARRA", "2010", 40000, 1000000)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $docId = $entry[1]
    $name = $entry[2]
    $period = $entry[3]
    $valorMora = $entry[4]
    $salarioBasico = $entry[5]

    $ws.Cells.Item($row, 3).Value = $docId
    $ws.Cells.Item($row, 4).Value = $name
    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
